$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy cell formatting (styles) for the two brand-new rows (215 and 216)
# from the existing last data row (214), for the columns that carry explicit styles (A and E).
$ws.Cells.Item(214, 1).Copy()
$ws.Cells.Item(215, 1).PasteSpecial(-4122)
$ws.Cells.Item(216, 1).PasteSpecial(-4122)
$ws.Cells.Item(214, 5).Copy()
$ws.Cells.Item(215, 5).PasteSpecial(-4122)
$ws.Cells.Item(216, 5).PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 91
$ws.Cells.Item(91, 1).Value = 89
$ws.Cells.Item(91, 2).Value = 6924568
$ws.Cells.Item(91, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(91, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(91, 5).Value = 45214.92013888889
$ws.Cells.Item(91, 6).Value = 'Atletico Morelia'
$ws.Cells.Item(91, 7).Value = 'Atlante'
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 1
$ws.Cells.Item(91, 10).Value = 'A'
$ws.Cells.Item(91, 11).Value = 2.4
$ws.Cells.Item(91, 12).Value = 3
$ws.Cells.Item(91, 13).Value = 2.875
$ws.Cells.Item(91, 14).Value = 2.7
$ws.Cells.Item(91, 15).Value = 3.1
$ws.Cells.Item(91, 16).Value = 2.8
$ws.Cells.Item(91, 17).Value = 0
$ws.Cells.Item(91, 18).Value = 1.85
$ws.Cells.Item(91, 19).Value = 1.95
$ws.Cells.Item(91, 20).Value = 2.25
$ws.Cells.Item(91, 21).Value = 1.975
$ws.Cells.Item(91, 22).Value = 1.725
$ws.Cells.Item(91, 23).Value = -1
$ws.Cells.Item(91, 24).Value = -1
$ws.Cells.Item(91, 25).Value = 1.8
$ws.Cells.Item(91, 26).Value = -1
$ws.Cells.Item(91, 27).Value = 0.95
$ws.Cells.Item(91, 28).Value = -1
$ws.Cells.Item(91, 29).Value = 0.7250000000000001

# Row 92
$ws.Cells.Item(92, 1).Value = 90
$ws.Cells.Item(92, 2).Value = 6924569
$ws.Cells.Item(92, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(92, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(92, 5).Value = 45214.92013888889
$ws.Cells.Item(92, 6).Value = 'Venados FC'
$ws.Cells.Item(92, 7).Value = 'Dorados'
$ws.Cells.Item(92, 8).Value = 4
$ws.Cells.Item(92, 9).Value = 1
$ws.Cells.Item(92, 10).Value = 'H'
$ws.Cells.Item(92, 11).Value = 1.615
$ws.Cells.Item(92, 12).Value = 4
$ws.Cells.Item(92, 13).Value = 4.5
$ws.Cells.Item(92, 14).Value = 1.5
$ws.Cells.Item(92, 15).Value = 4.75
$ws.Cells.Item(92, 16).Value = 5.75
$ws.Cells.Item(92, 17).Value = -1.25
$ws.Cells.Item(92, 18).Value = 1.925
$ws.Cells.Item(92, 19).Value = 1.875
$ws.Cells.Item(92, 20).Value = 3
$ws.Cells.Item(92, 21).Value = 1.75
$ws.Cells.Item(92, 22).Value = 1.95
$ws.Cells.Item(92, 23).Value = 0.5
$ws.Cells.Item(92, 24).Value = -1
$ws.Cells.Item(92, 25).Value = -1
$ws.Cells.Item(92, 26).Value = 0.925
$ws.Cells.Item(92, 27).Value = -1
$ws.Cells.Item(92, 28).Value = 0.75
$ws.Cells.Item(92, 29).Value = -1

# Row 210
$ws.Cells.Item(210, 1).Value = 208
$ws.Cells.Item(210, 2).Value = 7641711
$ws.Cells.Item(210, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(210, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(210, 5).Value = 45377.92013888889
$ws.Cells.Item(210, 6).Value = 'Cancun FC'
$ws.Cells.Item(210, 7).Value = 'Tapatio'
$ws.Cells.Item(210, 8).Value = 1
$ws.Cells.Item(210, 9).Value = 1
$ws.Cells.Item(210, 10).Value = 'D'
$ws.Cells.Item(210, 11).Value = 1.909
$ws.Cells.Item(210, 12).Value = 3.5
$ws.Cells.Item(210, 13).Value = 3.75
$ws.Cells.Item(210, 14).Value = 1.65
$ws.Cells.Item(210, 15).Value = 3.6
$ws.Cells.Item(210, 16).Value = 5
$ws.Cells.Item(210, 17).Value = -0.75
$ws.Cells.Item(210, 18).Value = 1.85
$ws.Cells.Item(210, 19).Value = 1.95
$ws.Cells.Item(210, 20).Value = 2.5
$ws.Cells.Item(210, 21).Value = 1.95
$ws.Cells.Item(210, 22).Value = 1.75
$ws.Cells.Item(210, 23).Value = -1
$ws.Cells.Item(210, 24).Value = 2.6
$ws.Cells.Item(210, 25).Value = -1
$ws.Cells.Item(210, 26).Value = -1
$ws.Cells.Item(210, 27).Value = 0.95
$ws.Cells.Item(210, 28).Value = -1
$ws.Cells.Item(210, 29).Value = 0.75

# Row 211
$ws.Cells.Item(211, 1).Value = 209
$ws.Cells.Item(211, 2).Value = 7641712
$ws.Cells.Item(211, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(211, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(211, 5).Value = 45378.00347222222
$ws.Cells.Item(211, 6).Value = 'Oaxaca'
$ws.Cells.Item(211, 7).Value = 'Tepatitlan FC'
$ws.Cells.Item(211, 8).Value = 1
$ws.Cells.Item(211, 9).Value = 0
$ws.Cells.Item(211, 10).Value = 'H'
$ws.Cells.Item(211, 11).Value = 1.833
$ws.Cells.Item(211, 12).Value = 3.6
$ws.Cells.Item(211, 13).Value = 3.75
$ws.Cells.Item(211, 14).Value = 1.909
$ws.Cells.Item(211, 15).Value = 3.5
$ws.Cells.Item(211, 16).Value = 4.2
$ws.Cells.Item(211, 17).Value = -0.5
$ws.Cells.Item(211, 18).Value = 1.825
$ws.Cells.Item(211, 19).Value = 1.975
$ws.Cells.Item(211, 20).Value = 2.25
$ws.Cells.Item(211, 21).Value = 1.85
$ws.Cells.Item(211, 22).Value = 1.95
$ws.Cells.Item(211, 23).Value = 0.909
$ws.Cells.Item(211, 24).Value = -1
$ws.Cells.Item(211, 25).Value = -1
$ws.Cells.Item(211, 26).Value = 0.825
$ws.Cells.Item(211, 27).Value = -1
$ws.Cells.Item(211, 28).Value = -1
$ws.Cells.Item(211, 29).Value = 0.95

# Row 212
$ws.Cells.Item(212, 1).Value = 210
$ws.Cells.Item(212, 2).Value = 7641713
$ws.Cells.Item(212, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(212, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(212, 5).Value = 45378.92013888889
$ws.Cells.Item(212, 6).Value = 'Universidad Guadalajara'
$ws.Cells.Item(212, 7).Value = 'Club Atletico La Paz'
$ws.Cells.Item(212, 8).Value = 2
$ws.Cells.Item(212, 9).Value = 4
$ws.Cells.Item(212, 10).Value = 'A'
$ws.Cells.Item(212, 11).Value = 1.666
$ws.Cells.Item(212, 12).Value = 3.8
$ws.Cells.Item(212, 13).Value = 4.333
$ws.Cells.Item(212, 14).Value = 1.444
$ws.Cells.Item(212, 15).Value = 4.5
$ws.Cells.Item(212, 16).Value = 7
$ws.Cells.Item(212, 17).Value = -1.25
$ws.Cells.Item(212, 18).Value = 1.925
$ws.Cells.Item(212, 19).Value = 1.875
$ws.Cells.Item(212, 20).Value = 2.75
$ws.Cells.Item(212, 21).Value = 1.9
$ws.Cells.Item(212, 22).Value = 1.9
$ws.Cells.Item(212, 23).Value = -1
$ws.Cells.Item(212, 24).Value = -1
$ws.Cells.Item(212, 25).Value = 6
$ws.Cells.Item(212, 26).Value = -1
$ws.Cells.Item(212, 27).Value = 0.875
$ws.Cells.Item(212, 28).Value = 0.8999999999999999
$ws.Cells.Item(212, 29).Value = -1

# Row 213
$ws.Cells.Item(213, 1).Value = 211
$ws.Cells.Item(213, 2).Value = 7641714
$ws.Cells.Item(213, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(213, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(213, 5).Value = 45379.00347222222
$ws.Cells.Item(213, 6).Value = 'Correcaminos'
$ws.Cells.Item(213, 7).Value = 'Cimarrones de Sonora FC'
$ws.Cells.Item(213, 8).Value = 0
$ws.Cells.Item(213, 9).Value = 2
$ws.Cells.Item(213, 10).Value = 'A'
$ws.Cells.Item(213, 11).Value = 1.95
$ws.Cells.Item(213, 12).Value = 3.25
$ws.Cells.Item(213, 13).Value = 3.6
$ws.Cells.Item(213, 14).Value = 2.2
$ws.Cells.Item(213, 15).Value = 3.25
$ws.Cells.Item(213, 16).Value = 3.5
$ws.Cells.Item(213, 17).Value = -0.25
$ws.Cells.Item(213, 18).Value = 1.85
$ws.Cells.Item(213, 19).Value = 1.95
$ws.Cells.Item(213, 20).Value = 2.25
$ws.Cells.Item(213, 21).Value = 1.925
$ws.Cells.Item(213, 22).Value = 1.875
$ws.Cells.Item(213, 23).Value = -1
$ws.Cells.Item(213, 24).Value = -1
$ws.Cells.Item(213, 25).Value = 2.5
$ws.Cells.Item(213, 26).Value = -1
$ws.Cells.Item(213, 27).Value = 0.95
$ws.Cells.Item(213, 28).Value = -0.5
$ws.Cells.Item(213, 29).Value = 0.4375

# Row 214
$ws.Cells.Item(214, 1).Value = 212
$ws.Cells.Item(214, 2).Value = 7641715
$ws.Cells.Item(214, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(214, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(214, 5).Value = 45379.92013888889
$ws.Cells.Item(214, 6).Value = 'Atletico Morelia'
$ws.Cells.Item(214, 7).Value = 'Venados FC'
$ws.Cells.Item(214, 8).Value = 3
$ws.Cells.Item(214, 9).Value = 4
$ws.Cells.Item(214, 10).Value = 'A'
$ws.Cells.Item(214, 11).Value = 2.6
$ws.Cells.Item(214, 12).Value = 3.4
$ws.Cells.Item(214, 13).Value = 2.4
$ws.Cells.Item(214, 14).Value = 2.45
$ws.Cells.Item(214, 15).Value = 3.2
$ws.Cells.Item(214, 16).Value = 3.1
$ws.Cells.Item(214, 17).Value = -0.25
$ws.Cells.Item(214, 18).Value = 2
$ws.Cells.Item(214, 19).Value = 1.8
$ws.Cells.Item(214, 20).Value = 2
$ws.Cells.Item(214, 21).Value = 1.9
$ws.Cells.Item(214, 22).Value = 1.9
$ws.Cells.Item(214, 23).Value = -1
$ws.Cells.Item(214, 24).Value = -1
$ws.Cells.Item(214, 25).Value = 2.1
$ws.Cells.Item(214, 26).Value = -1
$ws.Cells.Item(214, 27).Value = 0.8
$ws.Cells.Item(214, 28).Value = 0.8999999999999999
$ws.Cells.Item(214, 29).Value = -1

# Row 215
$ws.Cells.Item(215, 1).Value = 213
$ws.Cells.Item(215, 2).Value = 7640651
$ws.Cells.Item(215, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(215, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(215, 5).Value = 45380.00347222222
$ws.Cells.Item(215, 6).Value = 'Dorados'
$ws.Cells.Item(215, 7).Value = 'Atlante'
$ws.Cells.Item(215, 8).Value = 0
$ws.Cells.Item(215, 9).Value = 3
$ws.Cells.Item(215, 10).Value = 'A'
$ws.Cells.Item(215, 11).Value = 7.5
$ws.Cells.Item(215, 12).Value = 4.2
$ws.Cells.Item(215, 13).Value = 1.444
$ws.Cells.Item(215, 14).Value = 5.25
$ws.Cells.Item(215, 15).Value = 4
$ws.Cells.Item(215, 16).Value = 1.615
$ws.Cells.Item(215, 17).Value = 1
$ws.Cells.Item(215, 18).Value = 1.775
$ws.Cells.Item(215, 19).Value = 2.025
$ws.Cells.Item(215, 20).Value = 2.25
$ws.Cells.Item(215, 21).Value = 1.75
$ws.Cells.Item(215, 22).Value = 1.95
$ws.Cells.Item(215, 23).Value = -1
$ws.Cells.Item(215, 24).Value = -1
$ws.Cells.Item(215, 25).Value = 0.615
$ws.Cells.Item(215, 26).Value = -1
$ws.Cells.Item(215, 27).Value = 1.025
$ws.Cells.Item(215, 28).Value = 0.75
$ws.Cells.Item(215, 29).Value = -1

# Row 216
$ws.Cells.Item(216, 1).Value = 214
$ws.Cells.Item(216, 2).Value = 7641716
$ws.Cells.Item(216, 3).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(216, 4).Value = 'Mexico Liga de Expansion'
$ws.Cells.Item(216, 5).Value = 45380.875
$ws.Cells.Item(216, 6).Value = 'Club Celaya'
$ws.Cells.Item(216, 7).Value = 'Tlaxcala FC'
$ws.Cells.Item(216, 11).Value = 1.444
$ws.Cells.Item(216, 12).Value = 4.2
$ws.Cells.Item(216, 13).Value = 7.5
$ws.Cells.Item(216, 14).Value = 1.5
$ws.Cells.Item(216, 15).Value = 4.2
$ws.Cells.Item(216, 16).Value = 6.5
$ws.Cells.Item(216, 17).Value = -1
$ws.Cells.Item(216, 18).Value = 1.775
$ws.Cells.Item(216, 19).Value = 2.025
$ws.Cells.Item(216, 20).Value = 2.5
$ws.Cells.Item(216, 21).Value = 1.825
$ws.Cells.Item(216, 22).Value = 1.975
$ws.Cells.Item(216, 23).Value = 0
$ws.Cells.Item(216, 24).Value = 0
$ws.Cells.Item(216, 25).Value = 0
$ws.Cells.Item(216, 26).Value = 0
$ws.Cells.Item(216, 27).Value = 0

